$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.370.71"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.569.44"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3738"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.52"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07403"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.895"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "1.571.71"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001107"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06679"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.141"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "22.376.40"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.357"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.517"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.69%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "146.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.990"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "1.744.11"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.980"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9833"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.892"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.609"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08347"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.381"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02451"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2244"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06344"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.354"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6164"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.788"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5748"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.036"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.220"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("E51").Value = "  +0.46%  "
